# Re-process the data with the newly curated dimensions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update iaest-dimension:* labels to their curated equivalents.
$ws.Range("A2").Value = "iaest-measure:horas-trabajadas"
$ws.Range("C2").Value = "iaest-measure:ocupacion-1-digito-descripcion"
$ws.Range("E2").Value = "sdmx-dimension:refArea"

# Row 3: "dim" becomes "medida" for the now-measure columns A and C.
$ws.Range("A3").Value = "medida"
$ws.Range("C3").Value = "medida"

# Row 4: "skos:Concept" becomes "xsd:int" for measure columns, and
# "URI-Comunidad" for the refArea dimension column E.
$ws.Range("A4").Value = "xsd:int"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("E4").Value = "URI-Comunidad"

# Row 5 (mapping file references) is no longer needed; remove it entirely.
$ws.Range("A5:E5").EntireRow.Delete()
